# Generate Report for Handoff
# Updates the localization status report after a new XLIFF handoff
# generation run for the four "Ready for handoff" files
# (0db91efb, 3a04c8ba, 512d9a0f, 6f5a1543):
#  - Overview sheet: bump "Latest HO Xliff Generate Date" for those rows
#  - zh-cn / de-de sheets: bump "Priority" from "low" to "ht" and update
#    "Latest Handoff Datetime" for those rows

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-08-12 04:32:18"
$overview.Range("G5").Value = "2016-08-12 04:32:18"
$overview.Range("G6").Value = "2016-08-12 04:32:18"
$overview.Range("G7").Value = "2016-08-12 04:32:18"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-12 04:32:12"
$zhcn.Range("E5").Value = "ht"
$zhcn.Range("H5").Value = "2016-08-12 04:32:12"
$zhcn.Range("E6").Value = "ht"
$zhcn.Range("H6").Value = "2016-08-12 04:32:12"
$zhcn.Range("E7").Value = "ht"
$zhcn.Range("H7").Value = "2016-08-12 04:32:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "ht"
$dede.Range("H4").Value = "2016-08-12 04:32:18"
$dede.Range("E5").Value = "ht"
$dede.Range("H5").Value = "2016-08-12 04:32:18"
$dede.Range("E6").Value = "ht"
$dede.Range("H6").Value = "2016-08-12 04:32:18"
$dede.Range("E7").Value = "ht"
$dede.Range("H7").Value = "2016-08-12 04:32:18"
